# Vault / Observation Platform workbook update
# ----------------------------------------------
# Adds a new "types" lookup worksheet (type -> French long-name) in front of
# the existing data sheet, and extends the existing data sheet with a new
# "longname" column plus renames the "hydrophone" platform type to "mooring".

$wb = $excel.ActiveWorkbook

# The existing (only) sheet holds the observation-platform data; it keeps its
# underlying name ("Sheet1") but ends up as the second tab.
$data = $wb.Worksheets.Item(1)

# Insert a brand-new worksheet *before* the data sheet and name it "types".
$types = $wb.Worksheets.Add($data)
$types.Name = "types"

# --- Populate the "types" lookup table -------------------------------------
# Columns: A = numeric id, B = English type key, C = French long name.
# Cells are written in this particular order so that brand-new shared
# strings are interned in the same sequence as the source workbook
# (avion, bateau, mooring, mouillage, land, terre, space, espace,
#  underwater glider, planeur sous-marin, longname).

$types.Range("A1").Value = 1
$types.Range("B1").Value = "plane"
$types.Range("C1").Value = "avion"

$types.Range("A2").Value = 2
$types.Range("B2").Value = "boat"
$types.Range("C2").Value = "bateau"

$types.Range("A6").Value = 6
$types.Range("B6").Value = "mooring"
$types.Range("C6").Value = "mouillage"

$types.Range("A5").Value = 5
$types.Range("B5").Value = "land"
$types.Range("C5").Value = "terre"

$types.Range("A7").Value = 7
$types.Range("B7").Value = "space"
$types.Range("C7").Value = "espace"

$types.Range("A4").Value = 4
$types.Range("B4").Value = "underwater glider"
$types.Range("C4").Value = "planeur sous-marin"

$types.Range("A3").Value = 3
$types.Range("B3").Value = "drone"
$types.Range("C3").Value = "drone"

$types.Range("A1:C7").Select() | Out-Null

# --- Update the observation-platform data sheet -----------------------------
# Re-resolve by name since the sheet's position/collection entry shifted
# after the insert above.
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Row 13 used to describe a "hydrophone" platform; it is now a "mooring".
$sheet1.Range("A13").Value = "mooring"

# New trailing column: long (French) name, mirroring the "types" lookup.
$sheet1.Range("F1").Value = "longname"

$sheet1.Range("A21").Select() | Out-Null

# The data sheet remains the tab the user is looking at.
$sheet1.Activate() | Out-Null
